# DecidimLoginData.xlsx update:
#  - add a second worksheet "Full2" (a template/dupe of the login table, with
#    some rows repurposed for admin/user "template" logins and a block of
#    "INJECTED CODE" test rows)
#  - Full2 becomes the active sheet/tab
#  - Full1's old stale selection gets tidied up (tabSelected off, single-cell
#    selection on A3 instead of the old multi-row sqref)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- add the new sheet right after Full1 ------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Full2"

# ---- column widths / row heights for Full2 -----------------------------
# (25.3 is the COM "characters" width that lands closest to the ~26.24
# target once the host quantizes it to its internal pixel grid)
$ws2.Range("A1:B16").EntireRow.RowHeight = 28.7
$ws2.Columns.Item(1).ColumnWidth = 25.3
$ws2.Columns.Item(2).ColumnWidth = 25.3

# ---- header row ---------------------------------------------------------
$ws2.Range("A1").Value = "Decidim User Login"
$ws2.Range("B1").Value = "Decidim User Password"

# ---- row 2 / row 3: original admin + user logins, password "decidim" ---
$ws2.Range("A2").Value = "admin@example.org"
$ws2.Range("B2").Value = "decidim"
$ws2.Range("A3").Value = "user@example.org"
$ws2.Range("B3").Value = "decidim"

# ---- rows 4-6: userName3/4/5 template rows ------------------------------
$ws2.Range("A4").Value = "userName3"
$ws2.Range("B4").Value = "decidim123456"
$ws2.Range("A5").Value = "userName4"
$ws2.Range("B5").Value = "decidim123456"
$ws2.Range("A6").Value = "userName5"
$ws2.Range("B6").Value = "password5"

# ---- rows 7-10: ADMIN/USER template logins + repeated admin/user -------
$ws2.Range("A7").Value = "ADMIN@example.org"
$ws2.Range("B7").Value = "decidim123456"
$ws2.Range("A8").Value = "USER@example.org"
$ws2.Range("B8").Value = "decidim123456"
$ws2.Range("A9").Value = "admin@example.org"
$ws2.Range("B9").Value = "decidim123456"
$ws2.Range("A10").Value = "user@example.org"
$ws2.Range("B10").Value = "decidim123456"

# ---- rows 11-16: INJECTED CODE test rows --------------------------------
for ($r = 11; $r -le 16; $r++) {
    $ws2.Cells.Item($r, 1).Value = "INJECTED CODE"
    $ws2.Cells.Item($r, 2).Value = "INJECTED CODE"
}

# ---- hyperlinks on column A ---------------------------------------------
$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:admin@example.org", "", "", "admin@example.org") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "mailto:user@example.org", "", "", "user@example.org") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A7"), "mailto:ADMIN@example.org", "", "", "ADMIN@example.org") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A8"), "mailto:USER@example.org", "", "", "USER@example.org") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A9"), "mailto:admin@example.org", "", "", "admin@example.org") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A10"), "mailto:user@example.org", "", "", "user@example.org") | Out-Null

# match the blue hyperlink look already used on Full1's linked cells (Arial
# 10, blue font, no underline -- same as the xf Full1's A2/A3 already use)
foreach ($addr in @("A2", "A3", "A7", "A8", "A9", "A10")) {
    $rng = $ws2.Range($addr)
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 10
    $rng.Font.Color = 16711680
    $rng.Font.Underline = $false
}

# ---- selections / active tab --------------------------------------------
$ws1.Range("A3").Select() | Out-Null
$ws2.Range("A5").Select() | Out-Null
